$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete trailing rows 34 and 35 (data shrinks from F35 to F33)
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(34).Delete()

# Row 5: ('RM 14', -19.5, 12.3, -14.4, -5, 17.66) -> ('RM 14', -19.5, 12.3, -14.4, -5, None)
$ws.Range("F5").Value = ""

# Row 11: ('RM 58', -20.8, 11.4, -15.5, -7.9, None) -> ('RM 58', -20.8, 11.4, -15.5, -7.9, 17.65)
$ws.Range("F11").Value = 17.65

# Row 19: ('RM 125', -20.6, 13.2, None, -6.5, 17.81) -> ('RM 125', -20.6, 13.2, -15.5, -6.5, None)
$ws.Range("D19").Value = -15.5
$ws.Range("F19").Value = ""

# Row 21: ('RM 135', -18.9, None, -14.3, -8.7, 16.58) -> ('RM 135', -18.9, None, None, -8.7, 16.58)
$ws.Range("D21").Value = ""

# Row 23: ('RM 140', -19.5, None, None, -7, None) -> ('RM 140', -19.5, None, -13.9, -7, 16.48)
$ws.Range("D23").Value = -13.9
$ws.Range("F23").Value = 16.48

# Row 25: ('RM 145', -19.5, None, None, -7.1, None) -> ('RM 145', -19.5, None, None, -7.1, 16.6)
$ws.Range("F25").Value = 16.6

# Row 26: ('RM 232', -19.7, None, -15.6, -8.8, ' ') -> ('SC 5', -20.2, None, -13.8, -5, 17.38)
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

# Row 27: ('SC 5', -20.2, 10.8, -13.8, -5, 17.38) -> ('SC 101', -20.4, 10, None, -10, None)
$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = ""

# Row 28: ('SC 92', -17.2, 14.3, -14, -6.3, 17.22) -> ('SC 105', -19.6, 11.1, -13.7, -5.9, 17.44)
$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

# Row 29: ('SC 101', -20.4, None, -14.6, -10, 17) -> ('SC 119', -19.5, None, -13, -6.8, None)
$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = ""

# Row 30: ('SC 105', -19.6, 11.1, -13.7, -5.9, 17.44) -> ('SC 120', -19.7, 11.4, -13.6, -5.7, 16.89)
$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

# Row 31: ('SC 119', -19.5, 11.2, -13, -6.8, 18.06) -> ('SC 132', -18.8, 15.3, -13.7, -8.1, 17.18)
$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

# Row 32: ('SC 120', -19.7, 11.4, -13.6, -5.7, 16.89) -> ('SC 193', -19.9, 10.5, -14.7, -6.4, 17.39)
$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

# Row 33: ('SC 132', -18.8, 15.3, -13.7, -8.1, 17.18) -> ('SC 232', -19.5, 10.4, -14.1, -10.7, 17.53)
$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

